$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (column D) and volume-change (column E) figures for each
# coin row. Column D values frequently look like numbers/dates to Excel's
# auto-detection (e.g. "1.002", "291.32"), so we force the cell to Text
# format before writing, then restore the cell style afterwards so no
# stray number-format is left behind.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.443.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.575.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3771"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.95"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("E9").Value = "  +1.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.166"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07677"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.012"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.933"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.576.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06774"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.241"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("E23").Value = "  +1.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.442.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.754"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "145.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.046"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.749.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.248"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.015"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.014"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08574"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02564"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.35%  "
$ws.Range("E38").Value = "  +0.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.478"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6485"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.28%  "
$ws.Range("E44").Value = "  -2.62%  "
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6041"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.797"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.301"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.098"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("E50").Value = "  +3.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07336"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.01%  "
# Row 39/40: Hedera and TrustWalletToken swap places (with updated D/E values)
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06592"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.33%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.340"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.01%  "
